# Add the new "N=200000" worksheet at the end of the workbook, and populate
# it with the radix-sort timing results, matching the pattern used by the
# other "N=..." sheets already in the workbook.

$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "N=200000"

# Header row
$ws.Range("A1").Value = "Execução"
$ws.Range("B1").Value = "Tempo (ms)"

# Individual run measurements
$data = @(
    @(1, "534.7881 ms"),
    @(2, "534.4138 ms"),
    @(3, "537.4830 ms"),
    @(4, "535.7282 ms"),
    @(5, "615.5841 ms")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

# Summary rows
$ws.Range("A7").Value = "Média"
$ws.Range("B7").Value = "551.5995 ms"

$ws.Range("A8").Value = "Desvio Padrão"
$ws.Range("B8").Value = "35.7882 ms"
